$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = 44025
$ws.Range("B18").Value = 0.5
$ws.Range("C18").Value = "atrybuty w elementy"

$ws.Range("A19").Value = 44025
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = "Próba ogarnięcia XSLT. Ustawianie namespace. Czy jest sens w dynamicznym tworzeniu XSLT dla każdego templatu osobno?"

$ws.Range("A20").Value = 44025
$ws.Range("B20").Value = 0.5
$ws.Range("C20").Value = "Projektowanie rozwiązania problemu dynamicznego XSLT"

$ws.Range("A21").Value = 44025
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = "Poprawna nauka XSLT. Nowy plik output.xsl"

$null = $ws.Range("B21").Select()

